$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 14 de Mayo de 2020 a las 12:05"

$ws.Range("B15").Value = 78810
$ws.Range("C15").Value = 755
$ws.Range("D15").Value = 26675
$ws.Range("E15").Value = 49571
$ws.Range("G15").Value = 13
$ws.Range("H15").Value = 2564

$ws.Range("A37").Value = "Rumania"
$ws.Range("B37").Value = 16247
$ws.Range("C37").Value = 245
$ws.Range("D37").Value = 9053
$ws.Range("E37").Value = 6148
$ws.Range("F37").Value = 225
$ws.Range("G37").Value = 10
$ws.Range("H37").Value = 1046

$ws.Range("A38").Value = "Japon"
$ws.Range("B38").Value = 16049
$ws.Range("C38").Value = 0
$ws.Range("D38").Value = 8920
$ws.Range("E38").Value = 6451
$ws.Range("F38").Value = 243
$ws.Range("G38").Value = 0
$ws.Range("H38").Value = 678

$ws.Range("A39").Value = "Indonesia"
$ws.Range("B39").Value = 16006
$ws.Range("C39").Value = 568
$ws.Range("D39").Value = 3518
$ws.Range("E39").Value = 11445
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 15
$ws.Range("H39").Value = 1043

$ws.Range("B58").Value = 6145
$ws.Range("C58").Value = 91
$ws.Range("E58").Value = 1561

$ws.Range("B91").Value = 1464
$ws.Range("C91").Value = 1
$ws.Range("D91").Value = 267
$ws.Range("E91").Value = 1094
$ws.Range("F91").Value = 7

$ws.Range("B106").Value = 898
$ws.Range("C106").Value = 18
$ws.Range("D106").Value = 694
$ws.Range("E106").Value = 173

$ws.Range("B141").Value = 272
$ws.Range("C141").Value = 9
$ws.Range("E141").Value = 159

$ws.Range("A190").Value = "Butan"
$ws.Range("C190").Value = 8
$ws.Range("D190").Value = 5
$ws.Range("E190").Value = 14

$ws.Range("A191").Value = "Laos"
$ws.Range("B191").Value = 19
$ws.Range("E191").Value = 5

$ws.Range("A192").Value = "Fiyi"
$ws.Range("D192").Value = 14
$ws.Range("E192").Value = 4

$ws.Range("A193").Value = "Nueva Caledonia"
$ws.Range("D193").Value = 18
$ws.Range("H193").Value = 0

$ws.Range("A195").Value = "Belice"
$ws.Range("B195").Value = 18
$ws.Range("D195").Value = 16
$ws.Range("E195").Value = 0
$ws.Range("H195").Value = 2

$ws.Range("A196").Value = "Islas Virgenes de los Estados Unidos"
$ws.Range("D196").Value = 0
$ws.Range("E196").Value = 17

$ws.Range("A197").Value = "San Vicente y las Granadinas"
$ws.Range("B197").Value = 17
$ws.Range("C197").Value = 0
$ws.Range("D197").Value = 12
$ws.Range("E197").Value = 5
$ws.Range("H197").Value = 0

$ws.Range("A198").Value = "Mauritania"
$ws.Range("C198").Value = 1
$ws.Range("D198").Value = 6
$ws.Range("E198").Value = 8
$ws.Range("H198").Value = 2

$ws.Range("A199").Value = "Namibia"
$ws.Range("D199").Value = 12
$ws.Range("E199").Value = 4

$ws.Range("A200").Value = "Dominica"
$ws.Range("D200").Value = 15
$ws.Range("H200").Value = 0

$ws.Range("A201").Value = "Curazao"
$ws.Range("B201").Value = 16
$ws.Range("C201").Value = 0
$ws.Range("D201").Value = 14
$ws.Range("E201").Value = 1
$ws.Range("H201").Value = 1

$ws.Range("A215").Value = "Sahara Occidental"

$ws.Range("A216").Value = "Bonaire, San Eustaquio y Saba"
